$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Update B11: week 07-13/12/2015 -> 3+6+7
$ws.Range("B11").Formula = "=3+6+7"

# Update B12: week 14-20/12/2015 -> 5+3.5
$ws.Range("B12").Formula = "=5+3.5"

# Update selection to B11
$ws.Range("B11").Select()
